$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 76
$ws.Range("H76").Value = 2903.577
$ws.Range("I76").Value = 2927.8572
$ws.Range("J76").Value = 2801.6
$ws.Range("K76").Value = 2927.8572
$ws.Range("L76").Value = 2801.6
$ws.Range("M76").Value = -2612.8572
$ws.Range("N76").Value = -3431.6
# Row 79
$ws.Range("H79").Value = 2903.577
$ws.Range("I79").Value = 2927.8572
$ws.Range("J79").Value = 2801.6
$ws.Range("K79").Value = 2927.8572
$ws.Range("L79").Value = 2801.6
$ws.Range("M79").Value = -1835.8572
$ws.Range("N79").Value = -4985.6
# Row 125
$ws.Range("H125").Value = 1641.3
$ws.Range("I125").Value = 940
$ws.Range("J125").Value = 1816.625
$ws.Range("K125").Value = 8460
$ws.Range("L125").Value = 16349.625
$ws.Range("M125").Value = -6000
$ws.Range("N125").Value = -21269.625
# Row 131
$ws.Range("H131").Value = 3259.7273
$ws.Range("I131").Value = 1900
$ws.Range("J131").Value = 3395.7
$ws.Range("K131").Value = 5700
$ws.Range("L131").Value = 10187.1
$ws.Range("M131").Value = -660
$ws.Range("N131").Value = -20267.1
# Row 132
$ws.Range("H132").Value = 1573.3
$ws.Range("I132").Value = 1434.1666
$ws.Range("K132").Value = 4302.4998
$ws.Range("M132").Value = -1772.4998
# Row 135
$ws.Range("H135").Value = 1347.1428
$ws.Range("I135").Value = 503.8
$ws.Range("J135").Value = 6407.2
$ws.Range("K135").Value = 4534.2
$ws.Range("L135").Value = 57664.8
$ws.Range("M135").Value = -1999.2
$ws.Range("N135").Value = -62734.8
# Row 138
$ws.Range("H138").Value = 1892.8041
$ws.Range("I138").Value = 1677.9706
$ws.Range("J138").Value = 2008.746
$ws.Range("K138").Value = 5033.9118
$ws.Range("L138").Value = 6026.238
$ws.Range("M138").Value = 106.0882000000001
$ws.Range("N138").Value = -16306.238
# Row 141
$ws.Range("H141").Value = 1840.9822
$ws.Range("I141").Value = 740.2
$ws.Range("J141").Value = 11014.167
$ws.Range("K141").Value = 2220.6
$ws.Range("L141").Value = 33042.501
$ws.Range("M141").Value = 2959.4
$ws.Range("N141").Value = -43402.501

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 823815
$ws.Range("I32").Value = 962823.0600000001
$ws.Range("J32").Value = 17568.3
$ws.Range("K32").Value = 962823.0600000001
$ws.Range("L32").Value = 17568.3
$ws.Range("M32").Value = -962536.0600000001
$ws.Range("N32").Value = -18142.3
# Row 61
$ws.Range("H61").Value = 2244.5217
$ws.Range("I61").Value = 1781.6451
$ws.Range("J61").Value = 3201.1333
$ws.Range("K61").Value = 1781.6451
$ws.Range("L61").Value = 3201.1333
$ws.Range("M61").Value = -1569.6451
$ws.Range("N61").Value = -3625.1333
# Row 111
$ws.Range("H111").Value = 41111
$ws.Range("J111").Value = 41111
$ws.Range("L111").Value = 41111
$ws.Range("N111").Value = -49291
# Row 132
$ws.Range("H132").Value = 2972.1924
$ws.Range("I132").Value = 1968.238
$ws.Range("J132").Value = 7188.8
$ws.Range("K132").Value = 5904.714
$ws.Range("L132").Value = 21566.4
$ws.Range("M132").Value = -3374.714
$ws.Range("N132").Value = -26626.4
# Row 136
$ws.Range("H136").Value = 2244.5217
$ws.Range("I136").Value = 1781.6451
$ws.Range("J136").Value = 3201.1333
$ws.Range("K136").Value = 5344.9353
$ws.Range("L136").Value = 9603.3999
$ws.Range("M136").Value = -2794.9353
$ws.Range("N136").Value = -14703.3999

$ws = $wb.Worksheets.Item("BSM")
# Row 59
$ws.Range("H59").Value = 60780
$ws.Range("J59").Value = 60780
$ws.Range("L59").Value = 60780
$ws.Range("N59").Value = -62474
# Row 94
$ws.Range("H94").Value = 1171.52
$ws.Range("I94").Value = 1100.3529
$ws.Range("J94").Value = 1322.75
$ws.Range("K94").Value = 1100.3529
$ws.Range("L94").Value = 1322.75
$ws.Range("M94").Value = -649.3529000000001
$ws.Range("N94").Value = -2224.75
# Row 134
$ws.Range("H134").Value = 2178.2144
$ws.Range("I134").Value = 1687.1428
$ws.Range("J134").Value = 3651.4285
$ws.Range("K134").Value = 5061.428400000001
$ws.Range("L134").Value = 10954.2855
$ws.Range("M134").Value = -2526.428400000001
$ws.Range("N134").Value = -16024.2855

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3400.5166
$ws.Range("I31").Value = 891.01886
$ws.Range("J31").Value = 6900.6055
$ws.Range("K31").Value = 891.01886
$ws.Range("L31").Value = 6900.6055
$ws.Range("M31").Value = -596.01886
$ws.Range("N31").Value = -7490.6055
# Row 34
$ws.Range("H34").Value = 3400.5166
$ws.Range("I34").Value = 891.01886
$ws.Range("J34").Value = 6900.6055
$ws.Range("K34").Value = 891.01886
$ws.Range("L34").Value = 6900.6055
$ws.Range("M34").Value = -689.01886
$ws.Range("N34").Value = -7304.6055
# Row 58
$ws.Range("H58").Value = 1210.4147
$ws.Range("I58").Value = 968.04346
$ws.Range("J58").Value = 1520.1111
$ws.Range("K58").Value = 968.04346
$ws.Range("L58").Value = 1520.1111
$ws.Range("M58").Value = -765.04346
$ws.Range("N58").Value = -1926.1111
# Row 75
$ws.Range("H75").Value = 70000
$ws.Range("J75").Value = 70000
$ws.Range("L75").Value = 70000
$ws.Range("N75").Value = -71996
# Row 78
$ws.Range("H78").Value = 70000
$ws.Range("J78").Value = 70000
$ws.Range("L78").Value = 210000
$ws.Range("N78").Value = -219984
# Row 132
$ws.Range("H132").Value = 967.4761999999999
$ws.Range("I132").Value = 717.1667
$ws.Range("J132").Value = 2469.3333
$ws.Range("K132").Value = 2151.5001
$ws.Range("L132").Value = 7407.999899999999
$ws.Range("M132").Value = 378.4998999999998
$ws.Range("N132").Value = -12467.9999
# Row 134
$ws.Range("H134").Value = 3021.22
$ws.Range("I134").Value = 3252.2559
$ws.Range("J134").Value = 1602
$ws.Range("K134").Value = 9756.7677
$ws.Range("L134").Value = 4806
$ws.Range("M134").Value = -7221.7677
$ws.Range("N134").Value = -9876
# Row 136
$ws.Range("H136").Value = 1210.4147
$ws.Range("I136").Value = 968.04346
$ws.Range("J136").Value = 1520.1111
$ws.Range("K136").Value = 2904.13038
$ws.Range("L136").Value = 4560.3333
$ws.Range("M136").Value = -354.1303800000001
$ws.Range("N136").Value = -9660.3333

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 2075.2104
$ws.Range("I5").Value = 807.5
$ws.Range("J5").Value = 2413.2666
$ws.Range("K5").Value = 2422.5
$ws.Range("L5").Value = 7239.7998
$ws.Range("M5").Value = -2310.5
$ws.Range("N5").Value = -7463.7998
# Row 75
$ws.Range("H75").Value = 1000
$ws.Range("I75").Value = 1000
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 3000
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -2002
$ws.Range("N75").ClearContents()
# Row 78
$ws.Range("H78").Value = 1000
$ws.Range("I78").Value = 1000
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 9000
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -4008
$ws.Range("N78").ClearContents()
# Row 122
$ws.Range("H122").Value = 3569.8857
$ws.Range("I122").Value = 455.25
$ws.Range("K122").Value = 4097.25
$ws.Range("M122").Value = -1647.25
# Row 123
$ws.Range("H123").Value = 5618.5713
$ws.Range("I123").Value = 2332.5
$ws.Range("K123").Value = 6997.5
$ws.Range("M123").Value = -4547.5
# Row 126
$ws.Range("H126").Value = 3609.7144
$ws.Range("I126").Value = 1010
$ws.Range("J126").Value = 4318.727
$ws.Range("K126").Value = 3030
$ws.Range("L126").Value = 12956.181
$ws.Range("M126").Value = 1910
$ws.Range("N126").Value = -22836.181
# Row 131
$ws.Range("H131").Value = 3027.0876
$ws.Range("I131").Value = 456.4
$ws.Range("J131").Value = 3574.0425
$ws.Range("K131").Value = 1369.2
$ws.Range("L131").Value = 10722.1275
$ws.Range("M131").Value = 3670.8
$ws.Range("N131").Value = -20802.1275
# Row 135
$ws.Range("H135").Value = 2075.2104
$ws.Range("I135").Value = 807.5
$ws.Range("J135").Value = 2413.2666
$ws.Range("K135").Value = 7267.5
$ws.Range("L135").Value = 21719.3994
$ws.Range("M135").Value = -4732.5
$ws.Range("N135").Value = -26789.3994
# Row 138
$ws.Range("H138").Value = 4198.0557
$ws.Range("I138").Value = 1630
$ws.Range("J138").Value = 7408.125
$ws.Range("K138").Value = 4890
$ws.Range("L138").Value = 22224.375
$ws.Range("M138").Value = 250
$ws.Range("N138").Value = -32504.375
# Row 141
$ws.Range("H141").Value = 9832.895
$ws.Range("I141").Value = 9546.429
$ws.Range("J141").Value = 10000
$ws.Range("K141").Value = 28639.287
$ws.Range("L141").Value = 30000
$ws.Range("M141").Value = -23459.287
$ws.Range("N141").Value = -40360

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 3082.3489
$ws.Range("I132").Value = 2757.1794
$ws.Range("J132").Value = 6252.75
$ws.Range("K132").Value = 8271.538199999999
$ws.Range("L132").Value = 18758.25
$ws.Range("M132").Value = -5741.538199999999
$ws.Range("N132").Value = -23818.25

$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 2807.875
$ws.Range("I100").Value = 2810.5
$ws.Range("J100").Value = 2800
$ws.Range("K100").Value = 2810.5
$ws.Range("L100").Value = 2800
$ws.Range("M100").Value = -2269.5
$ws.Range("N100").Value = -3882
# Row 132
$ws.Range("H132").Value = 2871.4792
$ws.Range("I132").Value = 2731.6072
$ws.Range("K132").Value = 8194.821599999999
$ws.Range("M132").Value = -5664.821599999999
# Row 136
$ws.Range("H136").Value = 3877525.2
$ws.Range("I136").Value = 1689.2333
$ws.Range("J136").Value = 12821763
$ws.Range("K136").Value = 5067.699900000001
$ws.Range("L136").Value = 38465289
$ws.Range("M136").Value = -2517.699900000001
$ws.Range("N136").Value = -38470389

$ws = $wb.Worksheets.Item("WVR")
# Row 5
$ws.Range("H5").Value = 23692.77
$ws.Range("J5").Value = 23692.77
$ws.Range("L5").Value = 23692.77
$ws.Range("N5").Value = -23916.77
# Row 18
$ws.Range("H18").Value = 15750.875
$ws.Range("I18").Value = 8000
$ws.Range("K18").Value = 8000
$ws.Range("M18").Value = -7827
# Row 132
$ws.Range("H132").Value = 1465.6078
$ws.Range("I132").Value = 1399.45
$ws.Range("J132").Value = 1706.1818
$ws.Range("K132").Value = 4198.35
$ws.Range("L132").Value = 5118.5454
$ws.Range("M132").Value = -1668.35
$ws.Range("N132").Value = -10178.5454
# Row 136
$ws.Range("H136").Value = 1826.407
$ws.Range("I136").Value = 1574.7273
$ws.Range("J136").Value = 2656.95
$ws.Range("K136").Value = 4724.1819
$ws.Range("L136").Value = 7970.849999999999
$ws.Range("M136").Value = -2174.1819
$ws.Range("N136").Value = -13070.85

Write-Output "Edit applied successfully"
